$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-19: refreshed traffic incident data (destination sorting fix)
$data = @(
    @(2, "Roadwork", 1.324776067719318, 103.7420009258736, 566531145),
    @(3, "Vehicle breakdown", 1.335358592412027, 103.7750192680982, 139671206),
    @(4, "Roadwork", 1.33593426172264, 103.8099325412368, 392680115),
    @(5, "Roadwork", 1.294540495901838, 103.8761874369692, 366653031),
    @(6, "Roadwork", 1.337821500296478, 103.9805201027967, 4748287569),
    @(7, "Roadwork", 1.313207088444847, 103.9539123874291, 6690114550),
    @(8, "Roadwork", 1.324029203005592, 103.7433344949127, 566531128),
    @(9, "Roadwork", 1.294395580676698, 103.8665957683079, 26782016),
    @(10, "Roadwork", 1.295279710111475, 103.871392194302, 243494149),
    @(11, "Roadwork", 1.292970083024035, 103.8623491395152, 4732339357),
    @(12, "Roadwork", 1.426491921158762, 103.7821463731425, 1110672475),
    @(13, "Roadwork", 1.329408531652229, 103.846978150482, 245585883),
    @(14, "Roadwork", 1.395616215146429, 103.8576328064338, 5158040655),
    @(15, "Roadwork", 1.304490794689495, 103.770414048807, 247184656),
    @(16, "Roadwork", 1.408640772693746, 103.8080272748565, 148992674),
    @(17, "Roadwork", 1.394240890722329, 103.7742911599704, 158014654),
    @(18, "Roadwork", 1.340324439638299, 103.8051868390426, 139655727),
    @(19, "Roadwork", 1.352252843616421, 103.6957155183491, 567939126)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
    $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
    $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
    $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
}

Write-Host "Updated rows 2-19"
